# Regenerate merged AHB files
# 1. Rename header columns: "<Name>_old" -> "<Name>_FV2304" (left block, cols A-J)
#    and "<Name>_new" -> "<Name>_FV2310" (right block, cols L-U)
# 2. Turn the data range into an Excel Table (ListObject)
# 3. Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($baseNames[$i] + "_FV2304")
    $ws.Cells.Item(1, $i + 12).Value = ($baseNames[$i] + "_FV2310")
}

# Build the table over the full used range A1:U61
$tableRange = $ws.Range("A1:U61")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = "TableStyleMedium9"

# Freeze panes below row 1
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
